$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 9
$ws.Range("F3").Value = 11373
$ws.Range("F4").Value = 1287
$ws.Range("F5").Value = 1155
$ws.Range("F7").Value = 1221
$ws.Range("F8").Value = 166
$ws.Range("F9").Value = 950
$ws.Range("F11").Value = 2276
$ws.Range("G11").Value = 60
$ws.Range("F13").Value = 1115
$ws.Range("F14").Value = 871
$ws.Range("F15").Value = 576
$ws.Range("F16").Value = 857
$ws.Range("F17").Value = 1010
$ws.Range("F19").Value = 102
$ws.Range("F20").Value = 686
$ws.Range("F21").Value = 712
$ws.Range("F23").Value = 406
$ws.Range("F24").Value = 1057
$ws.Range("F25").Value = 63
$ws.Range("F26").Value = 449
$ws.Range("F27").Value = 535
$ws.Range("F28").Value = 195
$ws.Range("F29").Value = 268
$ws.Range("F30").Value = 269
$ws.Range("F31").Value = 635
$ws.Range("F32").Value = 2551
$ws.Range("F33").Value = 434
$ws.Range("F34").Value = 26
$ws.Range("F35").Value = 167
$ws.Range("F37").Value = 72
$ws.Range("F38").Value = 1508
$ws.Range("F39").Value = 428
$ws.Range("F40").Value = 133
$ws.Range("F41").Value = 65
$ws.Range("F42").Value = 109
$ws.Range("F47").Value = 64

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 93
$ws.Range("F16").Value = 102
$ws.Range("F20").Value = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2220
$ws.Range("F3").Value = 679
$ws.Range("F4").Value = 641

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2220
$ws.Range("F3").Value = 11373
$ws.Range("F4").Value = 641
$ws.Range("F5").Value = 1155
$ws.Range("F8").Value = 1221
$ws.Range("F10").Value = 166
$ws.Range("F12").Value = 950
$ws.Range("F13").Value = 2276
$ws.Range("G13").Value = 60
$ws.Range("F15").Value = 1115
$ws.Range("F16").Value = 871
$ws.Range("F17").Value = 576
$ws.Range("F18").Value = 857
$ws.Range("F19").Value = 1010
$ws.Range("F22").Value = 102
$ws.Range("F23").Value = 686
$ws.Range("F26").Value = 712
$ws.Range("F28").Value = 406
$ws.Range("F29").Value = 1057
$ws.Range("F30").Value = 63
$ws.Range("F31").Value = 449
$ws.Range("F32").Value = 535
$ws.Range("F33").Value = 195
$ws.Range("F34").Value = 268
$ws.Range("F35").Value = 2552
$ws.Range("F37").Value = 434
$ws.Range("F38").Value = 72
$ws.Range("F39").Value = 1508
$ws.Range("F40").Value = 428
$ws.Range("F41").Value = 133
$ws.Range("F42").Value = 65
$ws.Range("F44").Value = 109
$ws.Range("F48").Value = 64
